# Update countries & provincias Spain
# Applies the data refresh described by the commit/diff:
#  - Timestamp footer text updated (16:35 -> 17:05)
#  - Updated case counts for several countries (rows 4, 13, 19, 102)
#  - Mozambique inserted ahead of Birmania in the country list (rows 155-157 shift down,
#    with Mozambique receiving new figures and Birmania/Martinica's old figures shifting
#    one row down)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 17:05"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1689727
$ws.Range("C4").Value = 3291
$ws.Range("D4").Value = 451749
$ws.Range("E4").Value = 1138597

# --- Row 13: India ---
$ws.Range("B13").Value = 144069
$ws.Range("C13").Value = 5533
$ws.Range("E13").Value = 81225
$ws.Range("G13").Value = 93
$ws.Range("H13").Value = 4117

# --- Row 19: Chile ---
$ws.Range("B19").Value = 73997
$ws.Range("C19").Value = 4895
$ws.Range("D19").Value = 29302
$ws.Range("E19").Value = 43934
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 761

# --- Row 102: Kenia ---
$ws.Range("D102").Value = 402
$ws.Range("E102").Value = 832

# --- Rows 155-157: Mozambique inserted before Birmania, shifting
#     Birmania and Martinica down one row ---
# Row 155 becomes Mozambique with brand-new figures
$ws.Range("A155").Value = "Mozambique"
$ws.Range("B155").Value = 209
$ws.Range("C155").Value = 15
$ws.Range("D155").Value = 51
$ws.Range("E155").Value = 157
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 1

# Row 156 becomes Birmania carrying the figures that used to be on row 155
$ws.Range("A156").Value = "Birmania"
$ws.Range("B156").Value = 201
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 122
$ws.Range("E156").Value = 73
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 6

# Row 157 becomes Martinica carrying the figures that used to be on row 156
$ws.Range("A157").Value = "Martinica"
$ws.Range("B157").Value = 197
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 91
$ws.Range("E157").Value = 92
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 14

# Row 158 (Benin) stays the same - no change needed
